# Weekly update: insert a new record row at row 609 (pushing the
# existing rows 609:631 down to 610:632) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 609; existing rows 609-631 shift to 610-632.
$ws.Rows.Item(609).Insert()

# Populate the newly inserted row 609 with the new record's data.
$ws.Range("A609").Value = 8
$ws.Range("B609").Value = "Terminal La Palmera de La Serena"
$ws.Range("C609").Value = "Coquimbo"
$ws.Range("D609").Value = 44753
$ws.Range("E609").Value = 4
$ws.Range("F609").Value = 100112024
$ws.Range("G609").Value = "Choclo"
$ws.Range("H609").Value = "Dulce o Americano"
$ws.Range("I609").Value = "Primera"
$ws.Range("J609").Value = 400
$ws.Range("K609").Value = 45000
$ws.Range("L609").Value = 46000
$ws.Range("M609").Value = 45500
$ws.Range("N609").Value = "$/malla 70 unidades"
$ws.Range("O609").Value = "Región de Arica y Parinacota"
$ws.Range("P609").Value = 650
$ws.Range("Q609").Value = 70
$ws.Range("R609").Value = "Hortaliza"
